# Minor update to the results: refresh the "alignment_time" (column G)
# values for the bowtie2 rows with newer measurements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = "00:00:54"
$ws.Range("G5").Value  = "00:01:17"
$ws.Range("G8").Value  = "00:00:19"
$ws.Range("G11").Value = "00:05:17"
$ws.Range("G14").Value = "00:07:22"
$ws.Range("G17").Value = "00:01:45"
$ws.Range("G20").Value = "00:02:20"
$ws.Range("G23").Value = "00:04:45"
$ws.Range("G26").Value = "00:00:22"
$ws.Range("G29").Value = "00:00:55"
$ws.Range("G32").Value = "00:01:20"
$ws.Range("G35").Value = "00:00:17"
$ws.Range("G38").Value = "00:05:22"
$ws.Range("G41").Value = "00:07:25"
$ws.Range("G44").Value = "00:01:45"
$ws.Range("G47").Value = "00:02:22"
$ws.Range("G50").Value = "00:04:44"
$ws.Range("G53").Value = "00:00:24"
